$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions scheduled update).
# D-column "Price" cells are stored as plain text in the source sheet (values
# like "1.000" or "27.552.39" are not real numbers), so force Text format
# before assigning to stop Excel's automatic number coercion from dropping
# trailing zeros / reformatting them.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.552.39"
$ws.Range("E2").Value = "  -2.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.753.05"
$ws.Range("E3").Value = "  -3.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.35"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4485"
$ws.Range("E7").Value = "  +2.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3613"
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07508"
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.12"
$ws.Range("E10").Value = "  -6.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.105"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.70"
$ws.Range("E13").Value = "  -6.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.053"
$ws.Range("E14").Value = "  -4.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.177"
$ws.Range("E15").Value = "  -4.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.749.59"
$ws.Range("E16").Value = "  -3.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.86"
$ws.Range("E17").Value = "  -2.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001066"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06391"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.00"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.872"
$ws.Range("E22").Value = "  -5.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.602.65"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.104"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.73"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.46"
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.957.91"
$ws.Range("E28").Value = "  -3.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.121"
$ws.Range("E29").Value = "  -7.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.34"
$ws.Range("E30").Value = "  -3.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.085"
$ws.Range("E31").Value = "  -10.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09034"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.568"
$ws.Range("E33").Value = "  -7.39%  "
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.99"
$ws.Range("E35").Value = "  -8.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02310"
$ws.Range("E36").Value = "  -2.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05977"
$ws.Range("E39").Value = "  -3.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.962"
$ws.Range("E40").Value = "  -5.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.204"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.388"
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.800"
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.31"
$ws.Range("E45").Value = "  -4.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.715"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5882"
$ws.Range("E47").Value = "  -3.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.959"
$ws.Range("E48").Value = "  -3.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.61"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.162"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06868"
$ws.Range("E51").Value = "  -1.85%  "

# Rows 37/38: Algorand and TheSandbox swapped rank position
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6376"
$ws.Range("E37").Value = "  -3.60%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2088"
$ws.Range("E38").Value = "  -4.15%  "
